# This script reproduces the "Updated cryptos list" GitHub Actions commit:
# row-by-row Price (D) / Volume(1h) (E) refreshes, plus the Polygon / InternetComputer(DFINITY)
# row swap (rows 22-23, ranked #20/#21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text updates -----------------------------------------------------
# These new values are not parsed as numbers by Excel (they contain extra
# separators, percent signs/padding, or are plain names/URLs), so a simple
# Value assignment keeps them stored as text exactly like the original cells.
$ws.Range("D2").Value = "64.820.17"
$ws.Range("D3").Value = "3.145.57"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.146.06"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  -2.47%  "
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("E13").Value = "  +1.81%  "
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").Value = "3.659.93"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "64.916.64"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "3.146.76"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("E30").Value = "  +5.85%  "
$ws.Range("E31").Value = "  -1.92%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("E35").Value = "  -2.46%  "
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("E37").Value = "  +3.74%  "
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Value = "3.003.38"
$ws.Range("E42").Value = "  -3.71%  "
$ws.Range("E43").Value = "  -2.86%  "
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("E46").Value = "  -3.94%  "
$ws.Range("D47").Value = "0.0₃0576"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("E50").Value = "  -3.83%  "
$ws.Range("E51").Value = "  -0.87%  "

# --- Numeric-looking text updates --------------------------------------------
# These new values (e.g. "6.08", "28.00") WOULD be auto-converted to numbers by
# a direct Value assignment, but the source cells are text. Stage each value in
# a scratch cell formatted as Text, then Copy/PasteSpecial only the *values* into
# the destination so the destination keeps its original (unstyled) formatting and
# ends up holding the same text string as before.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "575.55"
$scratch.Copy()
$ws.Range("D5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "148.76"
$scratch.Copy()
$ws.Range("D6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "0.524"
$scratch.Copy()
$ws.Range("D9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "6.08"
$scratch.Copy()
$ws.Range("D11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "0.497"
$scratch.Copy()
$ws.Range("D12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "0.0000259"
$scratch.Copy()
$ws.Range("D13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "36.95"
$scratch.Copy()
$ws.Range("D14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "503.01"
$scratch.Copy()
$ws.Range("D20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "14.75"
$scratch.Copy()
$ws.Range("D21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "0.710"
$scratch.Copy()
$ws.Range("D22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "15.16"
$scratch.Copy()
$ws.Range("D23").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "7.68"
$scratch.Copy()
$ws.Range("D24").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "83.70"
$scratch.Copy()
$ws.Range("D25").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "2.81"
$scratch.Copy()
$ws.Range("D30").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "27.45"
$scratch.Copy()
$ws.Range("D31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "6.18"
$scratch.Copy()
$ws.Range("D34").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "6.44"
$scratch.Copy()
$ws.Range("D35").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "474.75"
$scratch.Copy()
$ws.Range("D38").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "2.96"
$scratch.Copy()
$ws.Range("D40").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "0.116"
$scratch.Copy()
$ws.Range("D43").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "28.00"
$scratch.Copy()
$ws.Range("D46").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "2.23"
$scratch.Copy()
$ws.Range("D50").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Value = "117.01"
$scratch.Copy()
$ws.Range("D51").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Clear()
$excel.CutCopyMode = 0
